$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the purchase route technical name to purchase_stock wherever it
#    appears in the route_ids column (Q2:Q28) - both the plain reference and
#    the comma-joined "...,stock.route_warehouse0_mto" variant pick up the
#    substring replacement.
$ws.Range("Q2:Q28").Replace("purchase.route_warehouse0_buy", "purchase_stock.route_warehouse0_buy") | Out-Null

# 2) Re-point the frozen-pane scroll position from Q11 towards K2 and adjust
#    the corresponding pane selections (mirrors scrolling the frozen sheet
#    back near the top-left of the data and selecting column Q).
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("K2").Select() | Out-Null
$win.FreezePanes = $true
$ws.Columns("Q:Q").Select() | Out-Null

# 3) Row 29 grows slightly taller.
$ws.Rows("29:29").RowHeight = 13.8

# 4) A new (blank, default-styled) cell shows up at Q29.
$ws.Range("Q29").ClearFormats()
